# "Doing Updates for Financials"
# Refresh the FEDU yearly financials figures (Income Statement, Balance
# Sheet and Cash Flow Statement sections) with the latest reported values.
# Only the numeric figures in columns D/E/F change; row 14 (Non Recurring,
# in the Operating Expenses block) no longer has numbers for the two most
# recent periods and is now reported as "NA", matching the existing "NA"
# label already used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Income Statement -------------------------------------------------
$ws.Range("D8").Value  = 44600
$ws.Range("E8").Value  = 30200
$ws.Range("F8").Value  = 13900

$ws.Range("D9").Value  = 16200
$ws.Range("E9").Value  = 12700
$ws.Range("F9").Value  = 8200

$ws.Range("D10").Value = 28400
$ws.Range("E10").Value = 17500
$ws.Range("F10").Value = 5800

$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"

$ws.Range("D17").Value = 35500
$ws.Range("E17").Value = 20800
$ws.Range("F17").Value = 12900

$ws.Range("D18").Value = 9100
$ws.Range("E18").Value = 9400

$ws.Range("D20").Value = 1000
$ws.Range("F20").Value = -4800

$ws.Range("D21").Value = 11100
$ws.Range("E21").Value = 5800
$ws.Range("F21").Value = -3800

$ws.Range("D23").Value = 10100
$ws.Range("E23").Value = 5500
$ws.Range("F23").Value = -3900

$ws.Range("D24").Value = 3900

$ws.Range("D26").Value = 6200
$ws.Range("E26").Value = 2600
$ws.Range("F26").Value = -4600

$ws.Range("D27").Value = 5400
$ws.Range("F27").Value = -4600

$ws.Range("D32").Value = -1000
$ws.Range("F32").Value = 4800

$ws.Range("D33").Value = 5400
$ws.Range("F33").Value = -4600

$ws.Range("D35").Value = 5400
$ws.Range("F35").Value = -4600

# --- Balance Sheet ------------------------------------------------------
$ws.Range("D41").Value = 86600
$ws.Range("E41").Value = 34300
$ws.Range("F41").Value = 6300

$ws.Range("E43").Value = 6700
$ws.Range("F43").Value = 6000

$ws.Range("D46").Value = 88300
$ws.Range("E46").Value = 41900
$ws.Range("F46").Value = 12700

$ws.Range("D47").Value = 23500

$ws.Range("D52").Value = 2200

$ws.Range("D54").Value = 117600
$ws.Range("E54").Value = 43900
$ws.Range("F54").Value = 13500

$ws.Range("D59").Value = 19900
$ws.Range("E59").Value = 18500
$ws.Range("F59").Value = 7800

$ws.Range("D60").Value = 19900
$ws.Range("E60").Value = 18500
$ws.Range("F60").Value = 7800

$ws.Range("F62").Value = 5900

$ws.Range("D66").Value = 20900
$ws.Range("E66").Value = 19400
$ws.Range("F66").Value = 13700

$ws.Range("E70").Value = 24300
$ws.Range("F70").Value = 3300

$ws.Range("F72").Value = -4500

$ws.Range("D76").Value = 96700
$ws.Range("F76").Value = -3500

# --- Cash Flow Statement --------------------------------------------------
$ws.Range("D81").Value = 5400
$ws.Range("F81").Value = -4600

$ws.Range("D89").Value = 14200
$ws.Range("E89").Value = 17700

$ws.Range("D94").Value = -25300

$ws.Range("D96").Value = -18100

$ws.Range("D100").Value = 68500
$ws.Range("E100").Value = 11100

$ws.Range("D101").Value = -5200
$ws.Range("E101").Value = 700

$ws.Range("D102").Value = 52300
$ws.Range("E102").Value = 28000
